# reduces size of project
# - removes the old row 10 (ser:144 / ser:140 / "live on Product Hunt" signin / ser:143 blog)
# - the former row 11 becomes the new row 10, with its content updated:
#     * the "course ser 3,4,5,6,7" entry is dropped
#     * the "blog ser 159" entry is dropped
#     * a new "blog ser 162" entry is added
#     * columns are re-packed left-to-right with no gaps

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 10 entirely; row 11 (already ht=255) shifts up to become row 10.
$ws.Rows.Item(10).Delete()

# Re-write row 10 contents/order to match the target layout.
$ws.Range("A10").Value = 43981
$ws.Range("B10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 162"
$ws.Range("C10").Value = "type: signin`nwidth: 2`nheight: 1`nh3: Raise / Sponsor Funds`np: Try out our new feature. Raise Funds for your next project or Join us in distributing rations.`nbutton.primary: Create a Ticket*goto(""/createticket"")`nbutton.secondary: View Tickets*goto(""/tickets"")`nsvg: /icons/stars.svg"
$ws.Range("D10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 161"
$ws.Range("E10").Value = "type: signin`nwidth: 2`nheight: 1`nh3.w-half: Sign up to get unlimited access to the entire content of zakatlists`nbutton.primary: Sign In*goto(""/signin/home"")`nbutton.secondary: Sign Up for Rs 300 / Month*goto(""/signup"")"
$ws.Range("F10").Value = "type: meetup`nwidth: 2`nheight: 1`nh3: Meetup coming in`ndate: 2020,6,5,10,30,0,0`nbutton.default: Speak*goto(""https://forms.gle/dyydXFRSsKzeH4hZ6"")`nbutton.default: Attend*goto(""https://youtu.be/vscn-HP932E"")`nbutton.default: Details*goto(""https://www.meetup.com/techshek/events/270179438/"")"
$ws.Range("G10").Value = "type: subscribe`nwidth: 2`nheight: 1`nh3: Subscribe to stay tuned to zakatlists`ninput: enter your email here`nbutton.default: Submit"
$ws.Range("H10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 160"
$ws.Range("I10").Value = "type: footer`nwidth: 6`nheight: 1`np.small: Eat from their fruits, and give the due alms on the day of harvest. <br> - Al Quran 6:141`nfacebook: https://facebook.com/zakatlists`ntwitter: https://twitter.com/zakatlists`nmakerlog: https://getmakerlog.com/@punch__lines "

# Remove any leftover cells beyond column I in row 10 (J10 held the old footer reference).
$ws.Range("J10:M10").Clear()

# Update the view so the new last row is in focus, matching the committed state.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("H10").Select() | Out-Null
